# Update column G ("K" = strikeouts, formerly "Strike#") values with the
# regenerated save_data for oviedo_luis.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 1
    6  = 1
    7  = 0
    8  = 3
    9  = 3
    10 = 0
    11 = 2
    12 = 4
    13 = 0
    14 = 2
    15 = 2
    16 = 0
    17 = 2
    18 = 1
    19 = 2
    20 = 1
    21 = 0
    22 = 4
    23 = 1
    24 = 1
    25 = 0
    26 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
